$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "https://tensorflow.blog/2023/10/03/book-roadmap/"

$ws.Range("D24").Value = "[근황] 논문 9편 억셉 & Open LLM 1등"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/223228054459"
